$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert new "Player Info" sheet before the existing sheet
# ---------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfoHeaderRange = $playerInfo.Range("A1:D1")
$playerInfoHeaderRange.Font.Bold = $true
$playerInfoHeaderRange.HorizontalAlignment = -4108
$playerInfoHeaderRange.VerticalAlignment = -4160
$playerInfoHeaderRange.Borders.LineStyle = 1

$playerInfo.Range("A1:D2").NumberFormat = "@"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
$col = 1
foreach ($h in $playerInfoHeaders) {
    $playerInfo.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$playerInfoRow2 = @("4425", "Babar Azam", "Right Handed", "Right Arm Off Break")
$col = 1
foreach ($v in $playerInfoRow2) {
    $playerInfo.Cells.Item(2, $col).Value = $v
    $col = $col + 1
}

# ---------------------------------------------------------------
# 2. Update "ODI Batting": header D1 + D-column MATCH_CODE values,
#    drop the two empty inlineStr cells at B37 / B46
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ODI Batting")
$ws.Range("D1").Value = "MATCH_CODE"
$ws.Range("D2:D96").NumberFormat = "@"

$matchCodes = @(
    "3803", "3814", "3819", "3858", "3859", "3861", "3863", "3879", "3883", "3921",
    "3925", "3926", "3928", "3930", "3932", "3939", "3943", "3944", "3972", "3973",
    "3975", "3977", "3981", "4017", "4018", "4019", "4034", "4037", "4044", "4045",
    "4050", "4079", "4081", "4082", "4084", "4087", "4103", "4104", "4105", "4110",
    "4114", "4172", "4174", "4176", "4177", "4178", "4194", "4197", "4200", "4201",
    "4204", "4223", "4225", "4227", "4237", "4238", "4241", "4244", "4247", "4287",
    "4292", "4294", "4297", "4300", "4304", "4308", "4319", "4324", "4334", "4337",
    "4340", "4349", "4375", "4376", "4432", "4433", "4434", "4458", "4459", "4460",
    "4472", "4473", "4476", "4564", "4565", "4567", "4586", "4590", "4592", "4634",
    "4638", "4641", "4686", "4688", "4690"
)

$r = 2
foreach ($code in $matchCodes) {
    $ws.Cells.Item($r, 4).Value = $code
    $r = $r + 1
}

$ws.Range("B37").Value = $null
$ws.Range("B46").Value = $null

# ---------------------------------------------------------------
# 3. Insert new "ODI Batting Extra" sheet after "ODI Batting"
# ---------------------------------------------------------------
$extra = $wb.Worksheets.Add($null, $ws)
$extra.Name = "ODI Batting Extra"

$extraHeaderRange = $extra.Range("A1:F1")
$extraHeaderRange.Font.Bold = $true
$extraHeaderRange.HorizontalAlignment = -4108
$extraHeaderRange.VerticalAlignment = -4160
$extraHeaderRange.Borders.LineStyle = 1

$extra.Range("A1:A21").NumberFormat = "@"
$extra.Range("C1:F21").NumberFormat = "@"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
$col = 1
foreach ($h in $extraHeaders) {
    $extra.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$extraRows = @(
    @("4433", $null, $null, $null, $null, "NO"),
    @("4434", 3, "13", "1", "44.96%", "NO"),
    @("4458", 3, "17", "0", "37.59%", "YES"),
    @("4459", 3, "6", "0", "9.57%", "NO"),
    @("4460", 3, "7", "3", "29.38%", "YES"),
    @("4472", 3, "0", "0", $null, "NO"),
    @("4473", 3, "4", "0", "9.74%", "NO"),
    @("4476", 3, "14", "4", "47.73%", "NO"),
    @("4564", 3, "6", "0", "25.33%", "NO"),
    @("4565", $null, $null, $null, $null, "NO"),
    @("4567", 3, "12", "0", "49.07%", "YES"),
    @("4586", $null, $null, $null, $null, "NO"),
    @("4590", $null, $null, $null, $null, "NO"),
    @("4592", 3, "0", "0", "0.37%", "NO"),
    @("4634", $null, $null, $null, $null, "NO"),
    @("4638", 3, "7", "0", "29.84%", "NO"),
    @("4641", 3, "7", "2", "44.17%", "NO"),
    @("4686", $null, $null, $null, $null, "NO"),
    @("4688", 3, "8", "1", "43.41%", "NO"),
    @("4690", $null, $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $col = 1
    foreach ($val in $row) {
        if ($val -ne $null) {
            $extra.Cells.Item($r, $col).Value = $val
        }
        $col = $col + 1
    }
    $r = $r + 1
}

